$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "view search results by different criteria" goes from done -> open,
#     clearing its date (F6) and tech (G6) values ---
$ws.Cells.Item(6,5).Value = "open"
$ws.Cells.Item(6,6).ClearContents()
$ws.Cells.Item(6,7).ClearContents()

# --- Insert a new "pagination" row before the old row 8 ("account") block ---
$ws.Rows("8:10").Insert()
$ws.Cells.Item(8,3).Value = "pagination"
$ws.Cells.Item(8,5).Value = "done"
$ws.Cells.Item(8,6).Value = 43988
$ws.Cells.Item(8,6).NumberFormat = "d-mmm"
$ws.Cells.Item(8,7).Value = "ajax "

# --- Insert 2 new rows before the old row 11 ("cart") block, now shifted to row 14 ---
$ws.Rows("14:15").Insert()
$ws.Cells.Item(15,3).Value = "cart items in sidebarB"
$ws.Cells.Item(15,5).Value = "done"
$ws.Cells.Item(15,6).Value = 43989
$ws.Cells.Item(15,6).NumberFormat = "d-mmm"
$ws.Cells.Item(14,2).Value = "browse"
$ws.Cells.Item(14,3).Value = "grid view/list view"

# --- Insert 2 new rows before the old row 14 ("checkout") block, now shifted to row 19 ---
$ws.Rows("19:20").Insert()
$ws.Cells.Item(19,3).Value = "cart totalprice in sidebarB"
$ws.Cells.Item(19,5).Value = "done"
$ws.Cells.Item(19,6).Value = 43989
$ws.Cells.Item(19,6).NumberFormat = "d-mmm"
$ws.Cells.Item(20,3).Value = "cart badge with quantity"
$ws.Cells.Item(20,5).Value = "done"
$ws.Cells.Item(20,6).Value = 43989
$ws.Cells.Item(20,6).NumberFormat = "d-mmm"

# --- Selection / active cell moves to E20:F20 ---
$ws.Range("E20:F20").Select()

# --- Add explicit page setup (orientation portrait) ---
$ws.PageSetup.Orientation = 1
